$wb = $excel.ActiveWorkbook

# Fix the Swiss sheet's selection to span the whole used range before
# we move away from it (mirrors the author clicking "Select All" there).
$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Range("A1:D12").Select()

# Duplicate the Swiss sheet (right after it) to use as the template for
# the new Portugal market sheet, then rename it.
$swiss.Copy($null, $swiss)
$portugal = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Name = "Portugal"

# Resize the columns to the narrower Portugal-specific widths.
$portugal.Columns.Item(1).ColumnWidth = 22.166666666666668
$portugal.Columns.Item(2).ColumnWidth = 18.5
$portugal.Columns.Item(3).ColumnWidth = 13.666666666666666
$portugal.Columns.Item(4).ColumnWidth = 12

# Update the market / user-story cell content (insert the new shared
# strings in the same order as the source workbook: code first, then
# market name).
$portugal.Range("B4").Value = "NGC-3479/T2404"
$portugal.Range("B2").Value = "Portugal Market"

# With the narrower columns the wrap-text cells in column D now take two
# lines, so the header/body rows grow to double height.
$portugal.Rows.Item(3).RowHeight = 28.8
$portugal.Rows.Item(4).RowHeight = 28.8
$portugal.Rows.Item(5).RowHeight = 28.8

# Leave the new sheet active with B2 selected.
$portugal.Range("B2").Select()
